$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.540.82"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "1.835.26"
$ws.Range("E3").Value = "  -0.80%  "

$ws.Range("D4").Value = "'1.008"
$ws.Range("E4").Value = "  +0.70%  "

$ws.Range("D5").Value = "'313.91"
$ws.Range("E5").Value = "  +0.21%  "

$ws.Range("D6").Value = "'1.008"
$ws.Range("E6").Value = "  +0.72%  "

$ws.Range("D7").Value = "'0.4239"
$ws.Range("E7").Value = "  +0.09%  "

$ws.Range("D8").Value = "'0.3664"
$ws.Range("E8").Value = "  +0.54%  "

$ws.Range("D9").Value = "'0.07239"
$ws.Range("E9").Value = "  -0.73%  "

$ws.Range("D10").Value = "'0.8654"
$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'20.71"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.964.07"
$ws.Range("E12").Value = "  +6.48%  "

$ws.Range("D13").Value = "'5.377"
$ws.Range("E13").Value = "  +0.63%  "

$ws.Range("D14").Value = "'6.500"
$ws.Range("E14").Value = "  -0.41%  "

$ws.Range("D15").Value = "'0.06976"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("D16").Value = "'1.010"
$ws.Range("E16").Value = "  +0.73%  "

$ws.Range("D17").Value = "'79.59"
$ws.Range("E17").Value = "  -0.18%  "

$ws.Range("D18").Value = "'0.000009005"
$ws.Range("E18").Value = "  +0.94%  "

$ws.Range("E19").Value = "  +0.88%  "

$ws.Range("D20").Value = "'15.45"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").Value = "27.655.99"
$ws.Range("E21").Value = "  -0.08%  "

$ws.Range("D22").Value = "'5.021"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("D23").Value = "'10.80"
$ws.Range("E23").Value = "  +3.81%  "

$ws.Range("D24").Value = "2.124.32"
$ws.Range("E24").Value = "  +2.68%  "

$ws.Range("D25").Value = "'1.969"
$ws.Range("E25").Value = "  -0.71%  "

$ws.Range("D26").Value = "'154.25"

$ws.Range("D27").Value = "'18.40"
$ws.Range("E27").Value = "  -3.03%  "

$ws.Range("D28").Value = "'5.237"
$ws.Range("E28").Value = "  -0.64%  "

$ws.Range("D29").Value = "'114.77"
$ws.Range("E29").Value = "  -6.08%  "

$ws.Range("D30").Value = "'1.822"
$ws.Range("E30").Value = "  -3.45%  "

$ws.Range("D31").Value = "'0.08882"
$ws.Range("E31").Value = "  +0.20%  "

$ws.Range("D32").Value = "'0.7692"
$ws.Range("E32").Value = "  +0.00%  "

$ws.Range("D33").Value = "'4.536"
$ws.Range("E33").Value = "  -0.49%  "

$ws.Range("D34").Value = "'2.957"
$ws.Range("E34").Value = "  -0.63%  "

$ws.Range("D35").Value = "'1.147"
$ws.Range("E35").Value = "  +3.39%  "

$ws.Range("D36").Value = "'1.008"
$ws.Range("E36").Value = "  +0.77%  "

$ws.Range("D37").Value = "'1.098"
$ws.Range("E37").Value = "  +0.19%  "

$ws.Range("D38").Value = "'0.05360"
$ws.Range("E38").Value = "  +0.06%  "

$ws.Range("D39").Value = "'0.01942"
$ws.Range("E39").Value = "  +0.56%  "

$ws.Range("D40").Value = "'2.825"
$ws.Range("E40").Value = "  +0.53%  "

$ws.Range("D41").Value = "'0.5106"
$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("D42").Value = "'0.1657"
$ws.Range("E42").Value = "  +0.40%  "

$ws.Range("D43").Value = "'6.766"
$ws.Range("E43").Value = "  -1.85%  "

$ws.Range("D44").Value = "'8.477"
$ws.Range("E44").Value = "  +2.00%  "

$ws.Range("D45").Value = "'10.48"
$ws.Range("E45").Value = "  +1.65%  "

$ws.Range("D46").Value = "'0.06531"
$ws.Range("E46").Value = "  -0.21%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4681"
$ws.Range("E47").Value = "  -1.84%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'104.99"
$ws.Range("E48").Value = "  -0.64%  "

$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "'1.008"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("D50").Value = "'1.617"
$ws.Range("E50").Value = "  -0.70%  "

$ws.Range("D51").Value = "'1.806"
$ws.Range("E51").Value = "  +4.99%  "
